# Lecture Note - Lecture 01
#
# On the "다음시간" (Next time) slide, the text run enumerating HTML tags
# currently reads "<div>,<span>,<table> ". Extend it to also mention the
# <a> tag: "<div>,<span>,<table>,<a> ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(41)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

$oldText = "<div>,<span>,<table> "
$newText = "<div>,<span>,<table>,<a> "

$fullText = $tr.Text
$idx = $fullText.IndexOf($oldText)
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, $oldText.Length)
    $target.Text = $newText
}
